$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New draw-result row appended by the nightly auto-update job.
$row = 96
$rng = $ws.Range("A" + $row + ":E" + $row)

# Temporarily force text storage so values round-trip exactly as plain
# strings (matching every other row in this sheet) instead of Excel
# auto-coercing "2025-12-21" / "251221" into a date serial / number.
$rng.NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2025-12-21"
$ws.Cells.Item($row, 2).Value = "Pick 3"
$ws.Cells.Item($row, 3).Value = "251221"
$ws.Cells.Item($row, 4).Value = "3-8-0"
$ws.Cells.Item($row, 5).Value = "2025-12-21T21:38:20.613+04:00"

# Drop the temporary number format again so the new row doesn't carry an
# explicit style index that the rest of the sheet doesn't have.
$rng.ClearFormats()
